$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This block of six "Times New Roman" bullet paragraphs (Ken.../space/
# I would still.../space/And do we know.../space) gets a brand-new bullet
# inserted at the front ("On the Schedule drop-down...") which pushes all
# the existing text down by one slot, and a new blank "space" bullet is
# appended at the end of the block to keep the same alternating pattern.
# ---------------------------------------------------------------------------

function StripTrailingMark($s) {
    return $s.Substring(0, $s.Length - 1)
}

$p3 = $d.Paragraphs.Item(3)
$p4 = $d.Paragraphs.Item(4)
$p5 = $d.Paragraphs.Item(5)
$p6 = $d.Paragraphs.Item(6)
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)

$textKen    = StripTrailingMark($p3.Range.Text)
$textSpace1 = StripTrailingMark($p4.Range.Text)
$textIWould = StripTrailingMark($p5.Range.Text)
$textSpace2 = StripTrailingMark($p6.Range.Text)
$textAnd    = StripTrailingMark($p7.Range.Text)
$textSpace3 = StripTrailingMark($p8.Range.Text)

$newFirstText = "On the Schedule drop-down, add an option for " + [char]0x201C + "Course" + [char]0x201D + " (or maybe " + [char]0x201C + "Course Offering" + [char]0x201D + ")."

# Shift the existing text down by one paragraph.
$p3.Range.Text = $newFirstText
$p4.Range.Text = $textKen
$p5.Range.Text = $textSpace1
$p6.Range.Text = $textIWould
$p7.Range.Text = $textSpace2
$p8.Range.Text = $textAnd

# Add a new trailing "space" bullet paragraph after the (now shifted)
# "And do we know..." paragraph, mirroring the existing spacer paragraphs.
$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = $textSpace3

# ---------------------------------------------------------------------------
# Move the (hidden) "_GoBack" bookmark from the end of the "To copy forward"
# paragraph to the end of the newly-typed first paragraph -- mirroring what
# Word itself does when new text is typed at that location.
#
# A collapsed Range built directly on the exact boundary right before a
# paragraph mark can't be used as the bookmark anchor directly, so a tiny
# unique marker is typed there, located with Find (which returns a Range
# that collapses cleanly), used to plant the bookmark, and then removed.
# ---------------------------------------------------------------------------

$p3 = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$marker = "@@BOOKMARKMARKER@@"
$insertPoint.InsertAfter($marker)

$markerRange = $d.Content
$markerRange.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$d.Content.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
